$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ARC_IC BM Prices")

# Update the "as of" date text in cell A3
$ws.Range("A3").Value = "April 10, 2025 1/"

# Update the "Projected (P) or Final (F) 2024 Actual ARC-IC Price" columns (K and N)
# Row 9 - Wheat
$ws.Range("K9").Value = 5.5
$ws.Range("N9").Value = 5.5

# Row 10 - Barley
$ws.Range("K10").Value = 6.5
$ws.Range("N10").Value = 6.5

# Row 11 - Oats
$ws.Range("K11").Value = 3.45
$ws.Range("N11").Value = 3.45

# Row 12 - Peanuts
$ws.Range("K12").Value = 0.255
$ws.Range("N12").Value = 0.255

# Row 13 - Corn
$ws.Range("K13").Value = 4.3499999999999996
$ws.Range("N13").Value = 4.3499999999999996

# Row 14 - Grain Sorghum
$ws.Range("K14").Value = 4.0999999999999996
$ws.Range("N14").Value = 4.0999999999999996

# Row 15 - Soybeans
$ws.Range("K15").Value = 9.9499999999999993
$ws.Range("N15").Value = 9.9499999999999993

# Row 17 - Lentils
$ws.Range("K17").Value = 0.34799999999999998
$ws.Range("N17").Value = 0.34799999999999998

# Row 18 - Canola
$ws.Range("K18").Value = 0.19900000000000001
$ws.Range("N18").Value = 0.19900000000000001

# Row 21 - Sunflower Seed
$ws.Range("K21").Value = 0.2155
$ws.Range("N21").Value = 0.2155

# Row 22 - Flaxseed
$ws.Range("K22").Value = 12.3
$ws.Range("N22").Value = 12.3

# Row 23 - Mustard Seed
$ws.Range("K23").Value = 0.48249999999999998
$ws.Range("N23").Value = 0.48249999999999998

# Row 24 - Rapeseed
$ws.Range("K24").Value = 0.2
$ws.Range("N24").Value = 0.2

# Row 25 - Safflower
$ws.Range("K25").Value = 0.30299999999999999
$ws.Range("N25").Value = 0.30299999999999999

# Row 26 - Crambe
$ws.Range("K26").Value = 0.24
$ws.Range("N26").Value = 0.24

# Row 28 - Seed Cotton 4/
$ws.Range("K28").Value = 0.33610000000000001
$ws.Range("N28").Value = 0.33610000000000001

# Row 29 - Rice (long grain)
$ws.Range("K29").Value = 0.14199999999999999
$ws.Range("N29").Value = 0.14199999999999999

# Row 30 - Rice (med/short grain) 5/
$ws.Range("K30").Value = 0.152
$ws.Range("N30").Value = 0.152

# Row 31 - Rice (temperate japonica)
$ws.Range("K31").Value = 0.22500000000000001
$ws.Range("N31").Value = 0.22500000000000001
